$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Eventos")

$ws.Range('K2').Value = 'https://s3.samsara.com/samsara-dashcam-videos/4006124/281474991109864/1750183612496/A5MTrC4t3W-camera-video-segment-driver-1750183614996.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSFI63V7VP%2F20250619%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250619T040005Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjELP%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJGMEQCIFl26XjLK7LvqewYAEmWPWlh1qZYEPBDLSsJmdpfrCoYAiBcVrUh5vS2QI8DsYCu42rj3lLwWvoZawzM5rrOd4DdtyrmAwic%2F%2F%2F%2F%2F%2F%2F%2F%2F%2F8BEAQaDDc4MTIwNDk0MjI0NCIMVqwGlny02XK6WgR4KroDd6Lauc7WxOT2qnaBp7BL1bEFp6DYp3rZ0pRyydkS7dShnic%2FLPPuhKPPsX2nMFmI%2Bjn2lbio9TTQEGpeerf7ZOLEwIxcFMixioZjg0IsO5FP1SMOKlUxzkng8psxM%2F49sByLEd5OfYTUoGmpXjHQdkqT7VeP3u%2FnSwUnx2WsomBr58zqAvPRULX3YUp%2Fs%2F7B79iHGWFTF9naJNFpa1wk%2FmqpJySjwapWqSUcadAzVi2T8%2FdFVY1o0eGck1nE6diVL0PXuz5ZeMZuMrn0E7w6Nm0mtlDd0hdqY0GbD9aIMXm3WIL71hEJNWwjnNTAZ2L4l8Nx8t6x0a4eXAuGhXdTgbKkaAQKwE%2B9Mfz3XIEUtl8OrUHOIRAE3hYMbpofON23wM0MuE%2FEJTq9CsewYg1YPOgP0kdw2lxxMNzz2tTC3JpJkqSnhxNu08Lpl7CXsyWNAWfJPYhI6D2pZMKoRX7JTkmtBI93XwtqcK1Go4Xz1NdxuuP8NN9886rVfH0g%2FtL67kOQUn16teKw6ILm2Ad011WRvt401q61KS2%2B2NHLm9gQdi%2F4j6mw4ULGeVFIpmp8FvluCRHvXHH%2BhTD2%2BM3CBjqmAVXyaXmC436nNUpEGOxB%2BS%2FLSOWmnxaRkj4SixrcxYBERqC4sBpMFDJack786jYzpsVrhSdFLOgUR%2B1h6KRL2B8IQJjSBMo4YhzIgkn0Gz7ql4HbfsirdoTOp9YjHEdOjOLCl13cF7lxlWP%2FbKt35Z875ylsYXYuHEHQHNGh7gHvFVOJk%2BUiTDgwd4Bqh8clhxIP9AV1yNgoKCNxhZ4v7i664wjWnXs%3D&X-Amz-SignedHeaders=host&response-expires=Thu%2C%2019%20Jun%202025%2012%3A00%3A05%20GMT&X-Amz-Signature=5db1f90b3c5753587a6ce9139b1a01657b8e28c85a35c60552c2d759e5b84cff'
$ws.Range('L2').Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991109864/1750183612496/TAdN8NokD3-camera-video-segment-1750183614996.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSFI63V7VP%2F20250619%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250619T040005Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjELP%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJGMEQCIFl26XjLK7LvqewYAEmWPWlh1qZYEPBDLSsJmdpfrCoYAiBcVrUh5vS2QI8DsYCu42rj3lLwWvoZawzM5rrOd4DdtyrmAwic%2F%2F%2F%2F%2F%2F%2F%2F%2F%2F8BEAQaDDc4MTIwNDk0MjI0NCIMVqwGlny02XK6WgR4KroDd6Lauc7WxOT2qnaBp7BL1bEFp6DYp3rZ0pRyydkS7dShnic%2FLPPuhKPPsX2nMFmI%2Bjn2lbio9TTQEGpeerf7ZOLEwIxcFMixioZjg0IsO5FP1SMOKlUxzkng8psxM%2F49sByLEd5OfYTUoGmpXjHQdkqT7VeP3u%2FnSwUnx2WsomBr58zqAvPRULX3YUp%2Fs%2F7B79iHGWFTF9naJNFpa1wk%2FmqpJySjwapWqSUcadAzVi2T8%2FdFVY1o0eGck1nE6diVL0PXuz5ZeMZuMrn0E7w6Nm0mtlDd0hdqY0GbD9aIMXm3WIL71hEJNWwjnNTAZ2L4l8Nx8t6x0a4eXAuGhXdTgbKkaAQKwE%2B9Mfz3XIEUtl8OrUHOIRAE3hYMbpofON23wM0MuE%2FEJTq9CsewYg1YPOgP0kdw2lxxMNzz2tTC3JpJkqSnhxNu08Lpl7CXsyWNAWfJPYhI6D2pZMKoRX7JTkmtBI93XwtqcK1Go4Xz1NdxuuP8NN9886rVfH0g%2FtL67kOQUn16teKw6ILm2Ad011WRvt401q61KS2%2B2NHLm9gQdi%2F4j6mw4ULGeVFIpmp8FvluCRHvXHH%2BhTD2%2BM3CBjqmAVXyaXmC436nNUpEGOxB%2BS%2FLSOWmnxaRkj4SixrcxYBERqC4sBpMFDJack786jYzpsVrhSdFLOgUR%2B1h6KRL2B8IQJjSBMo4YhzIgkn0Gz7ql4HbfsirdoTOp9YjHEdOjOLCl13cF7lxlWP%2FbKt35Z875ylsYXYuHEHQHNGh7gHvFVOJk%2BUiTDgwd4Bqh8clhxIP9AV1yNgoKCNxhZ4v7i664wjWnXs%3D&X-Amz-SignedHeaders=host&response-expires=Thu%2C%2019%20Jun%202025%2012%3A00%3A05%20GMT&X-Amz-Signature=202bf5f835faded895e14a6865ba94c01750fef633054770cae08ba1d3336376'
$ws.Range('K3').Value = 'https://s3.samsara.com/samsara-dashcam-videos/4006124/281474991152736/1750182961683/iBjUYEdqdu-camera-video-segment-driver-1750182964183.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSFI63V7VP%2F20250619%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250619T040005Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjELP%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJGMEQCIFl26XjLK7LvqewYAEmWPWlh1qZYEPBDLSsJmdpfrCoYAiBcVrUh5vS2QI8DsYCu42rj3lLwWvoZawzM5rrOd4DdtyrmAwic%2F%2F%2F%2F%2F%2F%2F%2F%2F%2F8BEAQaDDc4MTIwNDk0MjI0NCIMVqwGlny02XK6WgR4KroDd6Lauc7WxOT2qnaBp7BL1bEFp6DYp3rZ0pRyydkS7dShnic%2FLPPuhKPPsX2nMFmI%2Bjn2lbio9TTQEGpeerf7ZOLEwIxcFMixioZjg0IsO5FP1SMOKlUxzkng8psxM%2F49sByLEd5OfYTUoGmpXjHQdkqT7VeP3u%2FnSwUnx2WsomBr58zqAvPRULX3YUp%2Fs%2F7B79iHGWFTF9naJNFpa1wk%2FmqpJySjwapWqSUcadAzVi2T8%2FdFVY1o0eGck1nE6diVL0PXuz5ZeMZuMrn0E7w6Nm0mtlDd0hdqY0GbD9aIMXm3WIL71hEJNWwjnNTAZ2L4l8Nx8t6x0a4eXAuGhXdTgbKkaAQKwE%2B9Mfz3XIEUtl8OrUHOIRAE3hYMbpofON23wM0MuE%2FEJTq9CsewYg1YPOgP0kdw2lxxMNzz2tTC3JpJkqSnhxNu08Lpl7CXsyWNAWfJPYhI6D2pZMKoRX7JTkmtBI93XwtqcK1Go4Xz1NdxuuP8NN9886rVfH0g%2FtL67kOQUn16teKw6ILm2Ad011WRvt401q61KS2%2B2NHLm9gQdi%2F4j6mw4ULGeVFIpmp8FvluCRHvXHH%2BhTD2%2BM3CBjqmAVXyaXmC436nNUpEGOxB%2BS%2FLSOWmnxaRkj4SixrcxYBERqC4sBpMFDJack786jYzpsVrhSdFLOgUR%2B1h6KRL2B8IQJjSBMo4YhzIgkn0Gz7ql4HbfsirdoTOp9YjHEdOjOLCl13cF7lxlWP%2FbKt35Z875ylsYXYuHEHQHNGh7gHvFVOJk%2BUiTDgwd4Bqh8clhxIP9AV1yNgoKCNxhZ4v7i664wjWnXs%3D&X-Amz-SignedHeaders=host&response-expires=Thu%2C%2019%20Jun%202025%2012%3A00%3A05%20GMT&X-Amz-Signature=ffc7ece036befc08b6cf0df4c2cffa6d72d9a8600797bbab1f87946fdc92461e'
$ws.Range('L3').Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991152736/1750182961683/Foph7ZSROm-camera-video-segment-1750182964183.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSFI63V7VP%2F20250619%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250619T040005Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjELP%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJGMEQCIFl26XjLK7LvqewYAEmWPWlh1qZYEPBDLSsJmdpfrCoYAiBcVrUh5vS2QI8DsYCu42rj3lLwWvoZawzM5rrOd4DdtyrmAwic%2F%2F%2F%2F%2F%2F%2F%2F%2F%2F8BEAQaDDc4MTIwNDk0MjI0NCIMVqwGlny02XK6WgR4KroDd6Lauc7WxOT2qnaBp7BL1bEFp6DYp3rZ0pRyydkS7dShnic%2FLPPuhKPPsX2nMFmI%2Bjn2lbio9TTQEGpeerf7ZOLEwIxcFMixioZjg0IsO5FP1SMOKlUxzkng8psxM%2F49sByLEd5OfYTUoGmpXjHQdkqT7VeP3u%2FnSwUnx2WsomBr58zqAvPRULX3YUp%2Fs%2F7B79iHGWFTF9naJNFpa1wk%2FmqpJySjwapWqSUcadAzVi2T8%2FdFVY1o0eGck1nE6diVL0PXuz5ZeMZuMrn0E7w6Nm0mtlDd0hdqY0GbD9aIMXm3WIL71hEJNWwjnNTAZ2L4l8Nx8t6x0a4eXAuGhXdTgbKkaAQKwE%2B9Mfz3XIEUtl8OrUHOIRAE3hYMbpofON23wM0MuE%2FEJTq9CsewYg1YPOgP0kdw2lxxMNzz2tTC3JpJkqSnhxNu08Lpl7CXsyWNAWfJPYhI6D2pZMKoRX7JTkmtBI93XwtqcK1Go4Xz1NdxuuP8NN9886rVfH0g%2FtL67kOQUn16teKw6ILm2Ad011WRvt401q61KS2%2B2NHLm9gQdi%2F4j6mw4ULGeVFIpmp8FvluCRHvXHH%2BhTD2%2BM3CBjqmAVXyaXmC436nNUpEGOxB%2BS%2FLSOWmnxaRkj4SixrcxYBERqC4sBpMFDJack786jYzpsVrhSdFLOgUR%2B1h6KRL2B8IQJjSBMo4YhzIgkn0Gz7ql4HbfsirdoTOp9YjHEdOjOLCl13cF7lxlWP%2FbKt35Z875ylsYXYuHEHQHNGh7gHvFVOJk%2BUiTDgwd4Bqh8clhxIP9AV1yNgoKCNxhZ4v7i664wjWnXs%3D&X-Amz-SignedHeaders=host&response-expires=Thu%2C%2019%20Jun%202025%2012%3A00%3A05%20GMT&X-Amz-Signature=596ce3d3e5b737844056576703c48db7903e8481daebee1c17c0998c1aa05a5d'

